$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 21 (2020-11-02): extra hour worked, updated description ---
$ws.Range("C21").Value = 15
$ws.Range("E21").Value = "Planning out Raider.io API integration, filtering relevant character information, creating bottom appbar for better navigation options"

# --- Row 22 (2020-11-03): add description ---
$ws.Range("E22").Value = "Worked on R&D course"

# --- Row 23 (2020-11-04): add description ---
$ws.Range("E23").Value = "Worked on R&D course // Finished now"

# --- Row 24 (2020-11-05): hours worked + description ---
$ws.Range("B24").Value = 12
$ws.Range("C24").Value = 15
$ws.Range("E24").Value = "Created EU realm list, search page"

# --- Row 26 (2020-11-09): hours worked ---
$ws.Range("B26").Value = 10
$ws.Range("C26").Value = 15

# --- Highlight the weekly-total cells (column F) with a light accent fill ---
$totals = $ws.Range("F5,F10,F15,F20,F25,F30,F35")
$totals.Interior.ThemeColor = 5
$totals.Interior.TintAndShade = 0.79998168889431442

# --- Restore the active selection to where editing left off ---
$ws.Activate()
$ws.Range("E26").Select()
